# Auto-generated Excel COM-interop script
# Applies scheduled-runner market price / profit recalculation updates
# to the Sheets workbook, per sheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# Row 62
$ws.Range("H62").Value = 2820.3635
$ws.Range("J62").Value = 3554
$ws.Range("L62").Value = 3554
$ws.Range("N62").Value = -4802
# Row 65
$ws.Range("H65").Value = 2820.3635
$ws.Range("J65").Value = 3554
$ws.Range("L65").Value = 17770
$ws.Range("N65").Value = -24010
# Row 69
$ws.Range("H69").Value = 6013
$ws.Range("I69").Value = 6013
$ws.Range("K69").Value = 18039
$ws.Range("M69").Value = -17165
# Row 72
$ws.Range("H72").Value = 6013
$ws.Range("I72").Value = 6013
$ws.Range("K72").Value = 54117
$ws.Range("M72").Value = -49749
# Row 76
$ws.Range("H76").Value = 3240.4443
$ws.Range("I76").Value = 3132.5
$ws.Range("J76").Value = 4104
$ws.Range("K76").Value = 3132.5
$ws.Range("L76").Value = 4104
$ws.Range("M76").Value = -2817.5
$ws.Range("N76").Value = -4734
# Row 79
$ws.Range("H79").Value = 3240.4443
$ws.Range("I79").Value = 3132.5
$ws.Range("J79").Value = 4104
$ws.Range("K79").Value = 3132.5
$ws.Range("L79").Value = 4104
$ws.Range("M79").Value = -2040.5
$ws.Range("N79").Value = -6288
# Row 99
$ws.Range("H99").Value = 1609.4286
$ws.Range("J99").Value = 3199.8333
$ws.Range("L99").Value = 9599.499899999999
$ws.Range("N99").Value = -12595.4999
# Row 116
$ws.Range("H116").Value = 3355.303
$ws.Range("I116").Value = 2652.5789
$ws.Range("K116").Value = 2652.5789
$ws.Range("M116").Value = 789.4211
# Row 125
$ws.Range("H125").Value = 1805.2
$ws.Range("I125").Value = 1491.5
$ws.Range("J125").Value = 2118.9
$ws.Range("K125").Value = 13423.5
$ws.Range("L125").Value = 19070.1
$ws.Range("M125").Value = -10963.5
$ws.Range("N125").Value = -23990.1
# Row 132
$ws.Range("H132").Value = 3775730.8
$ws.Range("I132").Value = 4168728.2
$ws.Range("J132").Value = 2956
$ws.Range("K132").Value = 12506184.6
$ws.Range("L132").Value = 8868
$ws.Range("M132").Value = -12503654.6
$ws.Range("N132").Value = -13928
# Row 135
$ws.Range("H135").Value = 1371.5
$ws.Range("I135").Value = 1051.2632
$ws.Range("K135").Value = 9461.3688
$ws.Range("M135").Value = -6926.3688
# Row 137
$ws.Range("H137").Value = 2328222
$ws.Range("I137").Value = 2705212.2
$ws.Range("J137").Value = 3449.8333
$ws.Range("K137").Value = 8115636.600000001
$ws.Range("L137").Value = 10349.4999
$ws.Range("M137").Value = -8113086.600000001
$ws.Range("N137").Value = -15449.4999
# Row 138
$ws.Range("H138").Value = 2766.0881
$ws.Range("I138").Value = 1406.95
$ws.Range("J138").Value = 4707.7144
$ws.Range("K138").Value = 4220.85
$ws.Range("L138").Value = 14123.1432
$ws.Range("M138").Value = 919.1499999999996
$ws.Range("N138").Value = -24403.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2737.28
$ws.Range("I32").Value = 2431.753
$ws.Range("J32").Value = 5209.273
$ws.Range("K32").Value = 2431.753
$ws.Range("L32").Value = 5209.273
$ws.Range("M32").Value = -2144.753
$ws.Range("N32").Value = -5783.273
# Row 45
$ws.Range("H45").Value = 1491.9111
$ws.Range("I45").Value = 1009.6
$ws.Range("K45").Value = 1009.6
$ws.Range("M45").Value = -632.6
# Row 61
$ws.Range("H61").Value = 2232.1191
$ws.Range("I61").Value = 936.56525
$ws.Range("J61").Value = 3800.4211
$ws.Range("K61").Value = 936.56525
$ws.Range("L61").Value = 3800.4211
$ws.Range("M61").Value = -724.56525
$ws.Range("N61").Value = -4224.4211
# Row 74
$ws.Range("H74").Value = 700.5
$ws.Range("I74").Value = 684.3871
$ws.Range("K74").Value = 684.3871
$ws.Range("M74").Value = 189.6129
# Row 77
$ws.Range("H77").Value = 700.5
$ws.Range("I77").Value = 684.3871
$ws.Range("K77").Value = 3421.9355
$ws.Range("M77").Value = 946.0645
# Row 102
$ws.Range("H102").Value = 5079.9
$ws.Range("I102").Value = 4724.875
$ws.Range("J102").Value = 6500
$ws.Range("K102").Value = 4724.875
$ws.Range("L102").Value = 6500
$ws.Range("M102").Value = -3102.875
$ws.Range("N102").Value = -9744
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 136
$ws.Range("H136").Value = 2232.1191
$ws.Range("I136").Value = 936.56525
$ws.Range("J136").Value = 3800.4211
$ws.Range("K136").Value = 2809.69575
$ws.Range("L136").Value = 11401.2633
$ws.Range("M136").Value = -259.6957499999999
$ws.Range("N136").Value = -16501.2633

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2013.0869
$ws.Range("I105").Value = 1706.25
$ws.Range("J105").Value = 2714.4285
$ws.Range("K105").Value = 1706.25
$ws.Range("L105").Value = 2714.4285
$ws.Range("M105").Value = 40.75
$ws.Range("N105").Value = -6208.4285

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2275433.8
$ws.Range("I31").Value = 3127182.5
$ws.Range("J31").Value = 4104.1665
$ws.Range("K31").Value = 3127182.5
$ws.Range("L31").Value = 4104.1665
$ws.Range("M31").Value = -3126887.5
$ws.Range("N31").Value = -4694.1665
# Row 34
$ws.Range("H34").Value = 2275433.8
$ws.Range("I34").Value = 3127182.5
$ws.Range("J34").Value = 4104.1665
$ws.Range("K34").Value = 3127182.5
$ws.Range("L34").Value = 4104.1665
$ws.Range("M34").Value = -3126980.5
$ws.Range("N34").Value = -4508.1665
# Row 36
$ws.Range("H36").Value = 36701.168
$ws.Range("I36").Value = 3349.3333
$ws.Range("K36").Value = 3349.3333
$ws.Range("M36").Value = -2961.3333
# Row 40
$ws.Range("H40").Value = 36701.168
$ws.Range("I40").Value = 3349.3333
$ws.Range("K40").Value = 3349.3333
$ws.Range("M40").Value = -3189.3333
# Row 86
$ws.Range("H86").Value = 3363.6052
$ws.Range("I86").Value = 2659.5
$ws.Range("J86").Value = 4331.75
$ws.Range("K86").Value = 2659.5
$ws.Range("L86").Value = 4331.75
$ws.Range("M86").Value = -1536.5
$ws.Range("N86").Value = -6577.75
# Row 89
$ws.Range("H89").Value = 3363.6052
$ws.Range("I89").Value = 2659.5
$ws.Range("J89").Value = 4331.75
$ws.Range("K89").Value = 13297.5
$ws.Range("L89").Value = 21658.75
$ws.Range("M89").Value = -7681.5
$ws.Range("N89").Value = -32890.75
# Row 99
$ws.Range("H99").Value = 3057
$ws.Range("I99").Value = 1299.6666
$ws.Range("J99").Value = 4375
$ws.Range("K99").Value = 1299.6666
$ws.Range("L99").Value = 4375
$ws.Range("M99").Value = 198.3334
$ws.Range("N99").Value = -7371
# Row 105
$ws.Range("H105").Value = 2463.158
$ws.Range("I105").Value = 2480
$ws.Range("K105").Value = 2480
$ws.Range("M105").Value = -733
# Row 122
$ws.Range("H122").Value = 3151.6875
$ws.Range("I122").Value = 2820.3635
$ws.Range("J122").Value = 3880.6
$ws.Range("K122").Value = 8461.0905
$ws.Range("L122").Value = 11641.8
$ws.Range("M122").Value = -6011.0905
$ws.Range("N122").Value = -16541.8
# Row 126
$ws.Range("H126").Value = 3057
$ws.Range("I126").Value = 1299.6666
$ws.Range("J126").Value = 4375
$ws.Range("K126").Value = 3898.9998
$ws.Range("L126").Value = 13125
$ws.Range("M126").Value = -1428.9998
$ws.Range("N126").Value = -18065
# Row 132
$ws.Range("H132").Value = 2450.5642
$ws.Range("I132").Value = 1843.6923
$ws.Range("K132").Value = 5531.0769
$ws.Range("M132").Value = -3001.0769
# Row 134
$ws.Range("H134").Value = 1803.36
$ws.Range("I134").Value = 1028.5294
$ws.Range("J134").Value = 3449.875
$ws.Range("K134").Value = 3085.5882
$ws.Range("L134").Value = 10349.625
$ws.Range("M134").Value = -550.5881999999997
$ws.Range("N134").Value = -15419.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 645.2
$ws.Range("I5").Value = 375.21213
$ws.Range("J5").Value = 5100
$ws.Range("K5").Value = 1125.63639
$ws.Range("L5").Value = 15300
$ws.Range("M5").Value = -1013.63639
$ws.Range("N5").Value = -15524
# Row 92
$ws.Range("H92").Value = 2199.889
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 2412.375
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 7237.125
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -9733.125
# Row 122
$ws.Range("H122").Value = 1681.2727
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1832.6666
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 16493.9994
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -21393.9994
# Row 132
$ws.Range("H132").Value = 4476
$ws.Range("I132").Value = 2452
$ws.Range("K132").Value = 22068
$ws.Range("M132").Value = -19538
# Row 135
$ws.Range("H135").Value = 645.2
$ws.Range("I135").Value = 375.21213
$ws.Range("J135").Value = 5100
$ws.Range("K135").Value = 3376.90917
$ws.Range("L135").Value = 45900
$ws.Range("M135").Value = -841.9091699999999
$ws.Range("N135").Value = -50970

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3262.353
$ws.Range("I80").Value = 4220
$ws.Range("J80").Value = 2411.111
$ws.Range("K80").Value = 4220
$ws.Range("L80").Value = 2411.111
$ws.Range("M80").Value = -3222
$ws.Range("N80").Value = -4407.111
# Row 83
$ws.Range("H83").Value = 3262.353
$ws.Range("I83").Value = 4220
$ws.Range("J83").Value = 2411.111
$ws.Range("K83").Value = 21100
$ws.Range("L83").Value = 12055.555
$ws.Range("M83").Value = -16108
$ws.Range("N83").Value = -22039.555
# Row 97
$ws.Range("H97").Value = 3047.1428
$ws.Range("I97").Value = 2005
$ws.Range("J97").Value = 4436.6665
$ws.Range("K97").Value = 2005
$ws.Range("L97").Value = 4436.6665
$ws.Range("M97").Value = -1509
$ws.Range("N97").Value = -5428.6665
# Row 122
$ws.Range("H122").Value = 3305.5293
$ws.Range("I122").Value = 2162.8
$ws.Range("J122").Value = 3781.6667
$ws.Range("K122").Value = 6488.400000000001
$ws.Range("L122").Value = 11345.0001
$ws.Range("M122").Value = -4038.400000000001
$ws.Range("N122").Value = -16245.0001
# Row 132
$ws.Range("H132").Value = 3357.075
$ws.Range("I132").Value = 3044.182
$ws.Range("J132").Value = 3739.5
$ws.Range("K132").Value = 9132.545999999998
$ws.Range("L132").Value = 11218.5
$ws.Range("M132").Value = -6602.545999999998
$ws.Range("N132").Value = -16278.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 3451382.5
$ws.Range("I136").Value = 8336809.5
$ws.Range("J136").Value = 2845.8235
$ws.Range("K136").Value = 25010428.5
$ws.Range("L136").Value = 8537.470499999999
$ws.Range("M136").Value = -25007878.5
$ws.Range("N136").Value = -13637.4705

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 21455.6
$ws.Range("J96").Value = 50390
$ws.Range("L96").Value = 50390
$ws.Range("N96").Value = -53136
# Row 122
$ws.Range("H122").Value = 558247.4
$ws.Range("I122").Value = 668830.25
$ws.Range("K122").Value = 2006490.75
$ws.Range("M122").Value = -2006490.75

